$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-05-24 Friday" "2024-05-25 Saturday"

Replace-Text "970÷7=" "208÷8="
Replace-Text "535÷6=" "244÷2="
Replace-Text "432÷8=" "897÷3="
Replace-Text "256÷8=" "420÷7="
Replace-Text "984÷5=" "502÷4="
Replace-Text "675÷8=" "530÷2="
Replace-Text "328÷2=" "805÷9="
Replace-Text "346÷6=" "830÷2="
Replace-Text "196÷8=" "341÷4="
Replace-Text "917÷9=" "769÷3="
Replace-Text "394÷9=" "270÷6="
Replace-Text "524÷8=" "205÷4="
Replace-Text "266÷3=" "458÷3="
Replace-Text "628÷4=" "293÷3="
Replace-Text "783÷4=" "525÷8="
Replace-Text "125÷9=" "628÷5="
Replace-Text "739÷2=" "679÷8="
Replace-Text "670÷6=" "187÷7="
Replace-Text "917÷4=" "502÷7="
Replace-Text "708÷8=" "712÷6="
Replace-Text "969÷3=" "105÷7="
Replace-Text "497÷6=" "173÷6="
Replace-Text "701÷2=" "665÷6="
Replace-Text "212÷2=" "665÷4="
Replace-Text "857÷6=" "154÷6="
